# Apply the "No Detect Data" / "NA" corrections described in the commit:
# "Changed SummaryStats to only output maxes that are detected."
#
# Sheet "Alluvial for Mapping" (Max Cr / Max Date columns V & W) gets "No Detect Data"
# Sheet "Alluvial Exhibit" (Max Cr [ug/L] / Date of Max columns G & H) gets "NA"
# for the same set of wells (rows identified by their Location/Well ID).

$wb = $excel.ActiveWorkbook

$wsMapping = $wb.Worksheets.Item("Alluvial for Mapping")
$wsExhibit = $wb.Worksheets.Item("Alluvial Exhibit")

# Rows on "Alluvial for Mapping" sheet whose V/W (Max Cr / Max Date) values
# need to be replaced with "No Detect Data"
$mappingRows = @(18, 22, 24, 25, 29, 36, 38, 39)

foreach ($r in $mappingRows) {
    $wsMapping.Range("V$r").Value = "No Detect Data"
    $wsMapping.Range("W$r").Value = "No Detect Data"
}

# Rows on "Alluvial Exhibit" sheet whose G/H (Max Cr [ug/L] / Date of Max)
# values need to be replaced with "NA"
$exhibitRows = @(21, 25, 27, 28, 33, 40, 42, 43)

foreach ($r in $exhibitRows) {
    $wsExhibit.Range("G$r").Value = "NA"
    $wsExhibit.Range("H$r").Value = "NA"
}
